$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.643.22'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.525.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.91%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '202.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.44%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '553.33'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.97%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.513.05'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.609'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.70%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.656'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '63.49'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +12.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.143'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -7.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -7.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.90'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.096.30'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.531.94'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.93%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.48'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '67.489.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.88'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.33%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -5.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '394.11'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.09'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -8.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.01'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.84'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.26'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.89'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '720.59'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.17'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.13'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -13.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.77'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.25'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.68%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.77'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -9.42%  '

$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.399'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.08'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.92%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.082.77'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -4.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0685'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -13.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.56'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -11.97%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -11.38%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.49%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0413'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.08%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '138.52'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.30'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.89'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.57%  '
